$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update demand-center name (shared string reused across the sheet)
$ws.Range("A2").Value = "Lüderitz"

# Updated lat/lon and annual demand figures (demand now a literal value,
# no longer computed via a formula)
$ws.Range("B2").Value = -26.642877645011101
$ws.Range("C2").Value = 15.1439290700957
$ws.Range("D2").Value = 54000000

# The sheet previously carried a long tail of empty formatted rows
# (rows 3:38) below the single data row -- clear them out so the used
# range shrinks back down to just the header + one data row.
$ws.Range("A3:E38").Clear()

# Move the active selection to where the user left off after trimming
# the sheet.
$ws.Range("A9").Select()
